$d = $word.ActiveDocument

# --- Step 1: Replace the pseudocode placeholder paragraph with the full pseudocode ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*pseudocode here*") {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq -1) { throw "pseudocode placeholder paragraph not found" }
$pCell = $d.Paragraphs($targetIdx)
$rCell = $pCell.Range
$xmlPseudocode = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D9652D" w:rsidRPr="00970D86" w:rsidRDefault="00970D86" w:rsidP="00970D86"><w:pPr>              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">PROGRAM </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>checkIn</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">function </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>checkIn</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>():</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        PRINT “What type of pet do you have?”</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        INPUT the customer’s pet type</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        WHILE the customer''s pet type is not cat or dog</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            PRINT “Invalid pet type.  Please specify cat or dog.”</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            INPUT the customer’s pet type</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        ENDWHILE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        IF the customer''s pet type is a dog THEN</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            IF available dog space is greater than 0 THEN </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                RETURN true</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                RETURN false</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            IF available cat space is greater than 0 THEN </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                RETURN true</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                RETURN false</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        IF there is space available THEN</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            PRINT "Has this pet stayed with us before? (Y/N)"</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">         </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:t xml:space="preserve">   INPUT Y or N</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            WHILE input is not Y or N</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                PRINT "Invalid response.  Please enter Y or N."</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                INPUT Y or N</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ENDHILE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            IF customer is an existing client THEN</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                UPDATE existing pet information</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t xml:space="preserve">            ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                ADD new pet information</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            PRINT "How many days will the pet be staying with us?"</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            INPUT the number of days</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">           </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            IF the customer''s pet type is a dog THEN</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                IF the duration of the stay is greater than 2 days THEN</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                    SET grooming to TRUE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                    SET grooming to FALSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            IF the customer''s pet type is a dog</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                SET the dog''s space number</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                DECREMENT available dog space</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                PRINT dog''s space number</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                SET the cat''s space number</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                DECREMENT available cat space</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">                PRINT cat''s space number</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            ENDIF            </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        ELSE</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">            PRINT "No Vacancy."</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">        ENDIF</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:suppressAutoHyphens/>
            </w:pPr>
            <w:r>
              <w:t>END.</w:t>
            </w:r>
          </w:p>
        </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rCell.InsertXML($xmlPseudocode)

# --- Step 2: Fix up OOP Principles table (bookmark id shift, lastRenderedPageBreak, drop stray _GoBack) ---
$p1 = -1
$p4 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Inheritance*") { $p1 = $i }
    if ($t -like "Polymorphism*") { $p4 = $i }
}
if ($p1 -eq -1 -or $p4 -eq -1) { throw "OOP principles paragraphs not found" }
$rOop = $d.Range($d.Paragraphs($p1).Range.Start, $d.Paragraphs($p4).Range.End)
$xmlOop = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D9652D" w:rsidRDefault="00F16BBD" w:rsidP="00970D86"><w:pPr><w:suppressAutoHyphens/></w:pPr><w:bookmarkStart w:id="1" w:name="_gjdgxs" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="1"/><w:r><w:t>Inheritance</w:t></w:r></w:p><w:p w:rsidR="00F16BBD" w:rsidRDefault="00F16BBD" w:rsidP="00970D86"><w:pPr><w:suppressAutoHyphens/></w:pPr><w:r><w:t>Abstraction // Not really used here</w:t></w:r></w:p><w:p w:rsidR="00F16BBD" w:rsidRDefault="00F16BBD" w:rsidP="00970D86"><w:pPr><w:suppressAutoHyphens/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Encapsulation</w:t></w:r></w:p><w:p w:rsidR="00F16BBD" w:rsidRPr="00970D86" w:rsidRDefault="00F16BBD" w:rsidP="00970D86"><w:pPr><w:suppressAutoHyphens/></w:pPr><w:r><w:t>Polymorphism // Specifically overloading constructors</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rOop.InsertXML($xmlOop)

Write-Output "done"
